$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Step 1: Insert a new "Meta description" paragraph right after the first
#         paragraph (the H1 title "Play Cluedo Spinning Detectives Free |
#         Game Review").  The new paragraph has:
#           - a leading empty run
#           - a bold run containing "Meta description"
#           - a normal run containing the rest of the meta text
# -------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
# The newly created paragraph inherits the Heading1 style from paragraph 1;
# reset it back to Normal (body text) so it matches a plain paragraph.
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range

$boldLabel = "Meta description"
$restOfText = ": Learn about the gameplay structure, prizes, and symbols in Cluedo Spinning Detectives. Play for free and try your luck with the high volatility and impressive payouts."
$fullMetaText = $boldLabel + $restOfText

# Insert the whole run of text first (unformatted), then go back and bold
# only the "Meta description" label so we end up with two separate runs
# without corrupting the "current typing format" for later inserts.
$insertPoint = $d.Range($metaRange.Start, $metaRange.Start)
$insertPoint.InsertAfter($fullMetaText)

$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $boldLabel.Length)
$boldRange.Font.Bold = 1

# -------------------------------------------------------------------------
# Step 2: Remove the duplicated bold title paragraph
#         ("Play Cluedo Spinning Detectives Free | Game Review") that used
#         to sit right before the closing italic paragraph.
# -------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = $paraCount; $i -ge 2; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)
    if ($candidateText -eq "Play Cluedo Spinning Detectives Free | Game Review") {
        $candidate.Range.Delete()
        break
    }
}

# -------------------------------------------------------------------------
# Step 3: Replace the text of the final (italic) paragraph with the new
#         image-generation prompt, keeping the italic run formatting intact.
# -------------------------------------------------------------------------
$finalParaCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($finalParaCount)
$lastRange = $lastPara.Range
$lastTextLen = $lastRange.Text.TrimEnd([char]13, [char]7).Length
$lastContentRange = $d.Range($lastRange.Start, $lastRange.Start + $lastTextLen)

$newPromptText = "Please create a feature image in a cartoon style with a happy Maya warrior wearing glasses and incorporating elements of the Cluedo Spinning Detective game. You can include the Tudor Hall, magnifying glass, Stanze Bonus logo, interrogation point, and Clue logo in the image. The Maya warrior should be holding a Clue card or a magnifying glass to symbolize the game's detective aspect. The background could be a mix of blue and green colors to add some excitement and appeal to the image. With your creative skills, you can design an image that captures the essence of the game and attracts potential players. Thank you!"

$lastContentRange.Text = $newPromptText

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
